# Actualización automática 2025-12-03 13:30:07
#
# Refresh the "CUMPLIMIENTO MENSUAL" (monthly compliance) report with the
# latest PRESUPUESTO / VENTA / POR CUMPLIR figures per product group.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Re-apply header formatting (bold, centered, top-aligned, thin border) ---
$hdr = $ws3.Range("A1:F1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1         # xlContinuous

# --- Updated figures: PRESUPUESTO (C), VENTA (D), POR CUMPLIR (E) ---
# row 2 - 240X120 PORCELANATO
$ws3.Range("C2").Value = 136.08
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 136.08

# row 3 - 240X80 PORCELANATO
$ws3.Range("C3").Value = 2415
$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 2415

# row 4 - FREGADEROS DE COCINA
$ws3.Range("C4").Value = 199.5
$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 199.5

# row 5 - GRIFERIAS
$ws3.Range("C5").Value = 172.82
$ws3.Range("D5").Value = 0
$ws3.Range("E5").Value = 172.82

# row 6 - INODOROS
$ws3.Range("C6").Value = 920.4299999999999
$ws3.Range("D6").Value = 0
$ws3.Range("E6").Value = 920.4299999999999

# row 7 - LAVABOS
$ws3.Range("C7").Value = 263.97
$ws3.Range("D7").Value = 0
$ws3.Range("E7").Value = 263.97

# row 8 - NO RESURTIBLES (budget moved out, cells now blank)
$ws3.Range("C8").Value = $null
$ws3.Range("D8").Value = 0
$ws3.Range("E8").Value = $null

# row 9 - OTROS (budget moved in)
$ws3.Range("C9").Value = 350
$ws3.Range("D9").Value = 0
$ws3.Range("E9").Value = 350

# row 10 - PANELES DECORATIVOS
$ws3.Range("C10").Value = 407.52
$ws3.Range("D10").Value = 0
$ws3.Range("E10").Value = 407.52

# row 11 - PIEDRA SINTERIZADA
$ws3.Range("C11").Value = 1518.3
$ws3.Range("D11").Value = 0
$ws3.Range("E11").Value = 1518.3

# row 12 - PORCELANATO
$ws3.Range("C12").Value = 27207.6
$ws3.Range("D12").Value = 0
$ws3.Range("E12").Value = 27207.6

# row 13 - PUERTAS DE SEGURIDAD
$ws3.Range("C13").Value = 111.04
$ws3.Range("D13").Value = 0
$ws3.Range("E13").Value = 111.04

# row 14 - TOTAL
$ws3.Range("C14").Value = 33702.26
$ws3.Range("D14").Value = 0
$ws3.Range("E14").Value = 33702.26

# --- Leave the cursor parked on the detail sheet, then restore the
#     original active tab (VENTAS POR GRUPO) so the saved workbook still
#     opens on the first sheet, matching the source view state. ---
$ws3.Range("E18:E19").Select()
$ws1.Activate()
